$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 60 (shifts rows 60:199 down to 61:200,
# the former row 199 becomes row 200).
$ws.Rows.Item(60).Insert()

# Populate the freshly inserted row 60 with the new price record.
$ws.Cells.Item(60, 1).Value = 3
$ws.Cells.Item(60, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(60, 3).Value = "Coquimbo"
$ws.Cells.Item(60, 4).Value = 44498
$ws.Cells.Item(60, 5).Value = 5
$ws.Cells.Item(60, 6).Value = 100112039
$ws.Cells.Item(60, 7).Value = "Ciboulette"
$ws.Cells.Item(60, 8).Value = "Sin especificar"
$ws.Cells.Item(60, 9).Value = "Primera"
$ws.Cells.Item(60, 10).Value = 160
$ws.Cells.Item(60, 11).Value = 1500
$ws.Cells.Item(60, 12).Value = 1500
$ws.Cells.Item(60, 13).Value = 1500
$ws.Cells.Item(60, 14).Value = "$/docena de atados"
$ws.Cells.Item(60, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(60, 16).Value = 500
$ws.Cells.Item(60, 17).Value = 3
$ws.Cells.Item(60, 18).Value = "Hortaliza"
